$wb = $excel.ActiveWorkbook

$values = @{
    "C11" = 0.6069719124519111
    "D11" = 0.0
    "E11" = 0.3260270996830741
    "F11" = -0.039000000000000035
    "G11" = 1.5829618029997903
    "H11" = 16.12947350163202
    "I11" = 1.65808677867577
}

foreach ($sheet in $wb.Worksheets) {
    $ws = $sheet
    foreach ($addr in $values.Keys) {
        $ws.Range($addr).Value = $values[$addr]
    }
}
